$d = $word.ActiveDocument

# Locate the run that currently reads "进入项目子目录".
$findRng = $d.Content
$found = $findRng.Find.Execute("进入项目子目录", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Re-materialize a plain Range over the same span. (Calling InsertXML
    # directly on the Range returned by Find.Execute inserts *after* the
    # found text instead of replacing it; a freshly constructed Range does
    # not have that quirk.)
    $rng = $d.Range($findRng.Start, $findRng.End)

    # Replace the whole run with three runs:
    #   1) "进入项目"  - same formatting as the original run
    #   2) "根"        - same formatting, plus rFonts/w:hint="eastAsia"
    #   3) "目录"      - same formatting as the original run
    # Using InsertXML lets us control run boundaries and rPr precisely,
    # including the w:hint attribute that Word's object model does not
    # expose as a discrete Font property.
    $xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r w:rsidRPr="00412C66"><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="仿宋" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t>进入项目</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="仿宋" w:hAnsi="Courier New" w:cs="Courier New" w:hint="eastAsia"/></w:rPr><w:t>根</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="仿宋" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t>目录</w:t></w:r></w:p>'
    $rng.InsertXML($xmlFrag)
}
